# Applies the docx stat updates described by the commit.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Simple single-value cell replacements (rows by table position) ---
$t.Rows.Item(1).Cells.Item(1).Range.Text  = "0M"       # was 99.96
$t.Rows.Item(2).Cells.Item(1).Range.Text  = "0M"       # was 0.04
$t.Rows.Item(3).Cells.Item(1).Range.Text  = "0M"       # was 106
$t.Rows.Item(4).Cells.Item(1).Range.Text  = "302"      # was 122
$t.Rows.Item(5).Cells.Item(1).Range.Text  = "0.00001"  # was 0.00003
$t.Rows.Item(6).Cells.Item(1).Range.Text  = "0.00056"  # was 0.00039
$t.Rows.Item(9).Cells.Item(1).Range.Text  = "0.00019"  # was 0.00012
$t.Rows.Item(10).Cells.Item(1).Range.Text = "0.00021"  # was 0.00014
$t.Rows.Item(11).Cells.Item(1).Range.Text = "0.00023"  # was 0.00017
$t.Rows.Item(12).Cells.Item(1).Range.Text = "0.04078"  # was 0.01571

# --- Collapse the tab-separated multi-value rows down to a single value ---
$t.Rows.Item(44).Cells.Item(1).Range.Text = "99.96"
$t.Rows.Item(45).Cells.Item(1).Range.Text = "0.04"
$t.Rows.Item(46).Cells.Item(1).Range.Text = "106"
